$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.399.42'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '''1.838.14'
$ws.Range("E3").Value = '  -2.18%  '
$ws.Range("D4").Value = '''0.9999'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''260.12'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '''0.5228'
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("D8").Value = '''0.3237'
$ws.Range("E8").Value = '  -6.57%  '
$ws.Range("D9").Value = '''0.06776'
$ws.Range("E9").Value = '  -2.73%  '
$ws.Range("D10").Value = '''18.63'
$ws.Range("E10").Value = '  -7.67%  '
$ws.Range("D11").Value = '''0.7633'
$ws.Range("E11").Value = '  -5.50%  '
$ws.Range("D12").Value = '''0.07679'
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").Value = '''1.844.81'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '''88.47'
$ws.Range("E14").Value = '  -2.16%  '
$ws.Range("D15").Value = '''5.017'
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("D16").Value = '''1.000'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("E17").Value = '  -4.84%  '
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = '''0.000007902'
$ws.Range("E19").Value = '  -2.40%  '
$ws.Range("D20").Value = '''26.446.04'
$ws.Range("E20").Value = '  -2.01%  '
$ws.Range("D21").Value = '''2.076.76'
$ws.Range("E21").Value = '  -1.84%  '
$ws.Range("D22").Value = '''4.561'
$ws.Range("E22").Value = '  -4.11%  '
$ws.Range("D23").Value = '''9.443'
$ws.Range("E23").Value = '  -6.21%  '
$ws.Range("D24").Value = '''5.938'
$ws.Range("E24").Value = '  -4.28%  '
$ws.Range("D25").Value = '''144.65'
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("E26").Value = '  -5.78%  '
$ws.Range("D27").Value = '''1.655'
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("D28").Value = '''16.96'
$ws.Range("D29").Value = '''111.33'
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("D30").Value = '''4.181'
$ws.Range("E30").Value = '  -4.39%  '
$ws.Range("D31").Value = '''4.139'
$ws.Range("E31").Value = '  -4.62%  '
$ws.Range("D32").Value = '''0.08737'
$ws.Range("E32").Value = '  -1.98%  '
$ws.Range("D33").Value = '''0.04815'
$ws.Range("E33").Value = '  -2.67%  '
$ws.Range("D34").Value = '''1.121'
$ws.Range("E34").Value = '  -5.03%  '
$ws.Range("D35").Value = '''2.844'
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("D36").Value = '''0.7003'
$ws.Range("E36").Value = '  -5.25%  '
$ws.Range("D37").Value = '''3.062'
$ws.Range("E37").Value = '  -6.95%  '
$ws.Range("E38").Value = '  -4.98%  '
$ws.Range("D39").Value = '''2.180'
$ws.Range("E39").Value = '  -8.99%  '
$ws.Range("D40").Value = '''0.4823'
$ws.Range("E40").Value = '  -6.72%  '
$ws.Range("D41").Value = '''111.47'
$ws.Range("E41").Value = '  -4.42%  '
$ws.Range("D42").Value = '''0.8884'
$ws.Range("E42").Value = '  -7.68%  '
$ws.Range("D43").Value = '''6.088'
$ws.Range("E43").Value = '  -1.99%  '
$ws.Range("D44").Value = '''0.9996'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '''7.635'
$ws.Range("E45").Value = '  -6.27%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.4118'
$ws.Range("E46").Value = '  -8.86%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.05855'
$ws.Range("E47").Value = '  -1.76%  '
$ws.Range("D48").Value = '''8.985'
$ws.Range("E48").Value = '  -4.01%  '
$ws.Range("D49").Value = '''34.68'
$ws.Range("E49").Value = '  -4.46%  '
$ws.Range("E50").Value = '  -9.87%  '
$ws.Range("D51").Value = '''0.8802'
$ws.Range("E51").Value = '  -0.80%  '
